$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -20.292
$ws.Range("A6").Value = -22.291
$ws.Range("A7").Value = -19.898
$ws.Range("D7").Value = -8.154
$ws.Range("A8").Value = -22.223
$ws.Range("D11").Value = -7.056999999999999
$ws.Range("D12").Value = -7.181000000000002
$ws.Range("D15").Value = -8.196000000000002
$ws.Range("A16").Value = -21.879
$ws.Range("A20").Value = -20.009
$ws.Range("D20").Value = -7.714
$ws.Range("A21").Value = -20.013
$ws.Range("D21").Value = -8.115
$ws.Range("D22").Value = -7.970999999999999
$ws.Range("D23").Value = -7.997
$ws.Range("A28").Value = -22.046
$ws.Range("A29").Value = -21.343
$ws.Range("D29").Value = -7.51
$ws.Range("A30").Value = -21.658
$ws.Range("A32").Value = -21.648
$ws.Range("D34").Value = -7.904000000000001
$ws.Range("A40").Value = -19.869
$ws.Range("D42").Value = -7.952
$ws.Range("D43").Value = -7.943
$ws.Range("D44").Value = -7.831999999999999
$ws.Range("D45").Value = -7.525000000000001
$ws.Range("A46").Value = -21.801
$ws.Range("D46").Value = -8.420999999999999
$ws.Range("D50").Value = -8.046000000000001
$ws.Range("A51").Value = -22.162
$ws.Range("D51").Value = -8.301
$ws.Range("A52").Value = -22.333
$ws.Range("A57").Value = -22.567
$ws.Range("D57").Value = -8.214
$ws.Range("A59").Value = -22.354
$ws.Range("A62").Value = -22.111
$ws.Range("D65").Value = -7.784999999999999
$ws.Range("A66").Value = -21.551
$ws.Range("D66").Value = -7.557
$ws.Range("D67").Value = -7.203999999999999
$ws.Range("A73").Value = -20.082
$ws.Range("A74").Value = -21.244
$ws.Range("A77").Value = -20.342
$ws.Range("D79").Value = -7.6
$ws.Range("D84").Value = -8.300000000000001
$ws.Range("D87").Value = -8.022000000000002
$ws.Range("A92").Value = -21.64
$ws.Range("D92").Value = -6.654000000000001
$ws.Range("D97").Value = -8.484
$ws.Range("A100").Value = -22.217
